$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header suffixes to the FV2210 / FV2304 release
# tags (columns A1:J1 were the "_old" headers, L1:U1 the "_new" headers;
# K1 is the "diff" column and stays untouched).
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value2.Replace("_old", "_FV2210")
}

for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value2.Replace("_new", "_FV2304")
}

# Turn the used range into a real Excel Table (ListObject) bound to the
# renamed headers.
$rng = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
